# Fill in the missing outlier-distance (C) and group-size (D) values
# for rows 8 and 9, matching the pattern already present in rows 2-7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 4.9730000495910645
$ws.Range("D8").Value = 20.0
$ws.Range("C9").Value = 1.8899999856948853
$ws.Range("D9").Value = 20.0
